$d = $word.ActiveDocument

# Find the paragraph whose run text is "Visualizing hierarchies" (the document's
# first heading-style paragraph) and bump its font size to 20pt (40 half-points),
# matching both the paragraph mark run properties and the run itself.
foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text.TrimEnd([char]13, [char]7) -eq "Visualizing hierarchies") {
        $r.Font.Size = 20
        $r.Font.SizeBi = 20
        break
    }
}
